# Move all renewables to guaranteed dispatch and out of least cost dispatch
# with new logit dispatch function.
#
# On the "BGDPbES" sheet, the BAU Guaranteed Dispatch Percentage (column B,
# year 2015) for each renewable electricity source is changed from 0 to 1
# (i.e. fully guaranteed dispatch). The remaining year columns (C:AK) hold
# formulas that simply reference column B for the same row, so they pick up
# the new value automatically on recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BGDPbES")

# Row 6  -> onshore wind
# Row 7  -> solar PV
# Row 8  -> solar thermal
# Row 14 -> offshore wind
$ws.Range("B6").Value = 1
$ws.Range("B7").Value = 1
$ws.Range("B8").Value = 1
$ws.Range("B14").Value = 1

$excel.Calculate()

# Leave the cursor/selection on BGDPbES where the editor last left it,
# then switch back to the About sheet (the sheet that remains active
# when the workbook is saved).
$ws.Activate()
$ws.Range("B15").Select()

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
$wsAbout.Range("C45").Select()
